$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 ("Docentes responsáveis:" row gains its
# own value row below it, like "Objetivos:"/"Objectives:" above), shifting
# everything from the old row 13 onward down by one.
$ws.Rows.Item(13).Insert()

# --- Objetivos / Objectives (PT objective text replaces the placeholder
#     professor name that had been miswired into row 10) ---
$ws.Range("B10").Value = 'Fornecer aos alunos os conceitos básicos e técnicas de dimensionamento dos principais processos e operações unitárias envolvidas no escoamento de fluidos, sistemas particulados e troca térmica.'
$ws.Range("C10").Value = 'Fornecer aos alunos os conceitos básicos e técnicas de dimensionamento dos principais processos e operações unitárias envolvidas no escoamento de fluidos, sistemas particulados e troca térmica.'

# --- Docentes responsáveis: (new row 13 gets the professor name) ---
$ws.Range("B13").Value = '4780627 - Ana Lucia Gabas Ferreira'
$ws.Range("C13").Value = '4780627 - Ana Lucia Gabas Ferreira'

# --- Programa resumido: (Portuguese short syllabus replaces "Semestral") ---
$ws.Range("B14").Value = 'Operações unitárias e processos: reologia de fluidos, dimensionamento de tubulações e acessórios, bombeamento, agitação e mistura, caracterização de partículas e leito de partículas, sedimentação, filtração, processos com membranas. Operações unitárias de troca térmica: trocadores de calor e evaporadores.'
$ws.Range("C14").Value = 'Operações unitárias e processos: reologia de fluidos, dimensionamento de tubulações e acessórios, bombeamento, agitação e mistura, caracterização de partículas e leito de partículas, sedimentação, filtração, processos com membranas. Operações unitárias de troca térmica: trocadores de calor e evaporadores.'

# --- Programa: (Portuguese full syllabus replaces the stray date value) ---
$ws.Range("B16").Value = '- Reologia de fluidos,- Dimensionamento de tubulações,- Acessórios e bombeamento para fluidos industriais,- Agitação e mistura,- Caracterização de partículas e leito de partículas,- Sedimentação,- Filtração,- Processos com membranas.- Operações unitárias de troca térmica: trocadores de calor e evaporadores.'
$ws.Range("C16").Value = '- Reologia de fluidos,- Dimensionamento de tubulações,- Acessórios e bombeamento para fluidos industriais,- Agitação e mistura,- Caracterização de partículas e leito de partículas,- Sedimentação,- Filtração,- Processos com membranas.- Operações unitárias de troca térmica: trocadores de calor e evaporadores.'

# --- Método: (evaluation method text, shifted up from the old Critério: row) ---
$ws.Range("B19").Value = 'Avaliação composta por duas provas.'
$ws.Range("C19").Value = 'Avaliação composta por duas provas.'

# --- Critério: (grading criteria text) ---
$ws.Range("B20").Value = 'Média das notas das provas.'
$ws.Range("C20").Value = 'Média das notas das provas.'

# --- Norma de recuperação: (makeup exam rule text) ---
$ws.Range("B21").Value = 'Prova única com todo o conteúdo da disciplina, sendo que a nota [(nota final do semestre + nota de recuperação)/2] deverá ser igual ou superior a 5,0 (cinco).'
$ws.Range("C21").Value = 'Prova única com todo o conteúdo da disciplina, sendo que a nota [(nota final do semestre + nota de recuperação)/2] deverá ser igual ou superior a 5,0 (cinco).'

# --- Bibliografia: (actual bibliography content, previously absent) ---
$ws.Range("B22").Value = 'Bibliografia básica:DI BERNARDO, L., Métodos e Técnicas de Tratamento de Água, ABES, Rio de Janeiro, Brasil, 1992.FOUST, A.S., WENZEL, L. A., CLUMP, C.W., MAUS, L., ANDERSEN, L.B. Princípio das operações unitárias. Rio de Janeiro: Editora Guanabara Dois, 1982.GEANKOPLIS, C.J. Procesos de transporte y operaciones unitarias. Compañía Editorial Continental, S.A. de C.V. México, D.F., 1998.PERRY, R.H. and CHILTON, C.H. Manual de Engenharia Química. 5a ed., Guanabara Dois, Rio de Janeiro, 1986.REYNOLDS, T.D.; RICHARDS, P. Unit Operations and Processes in environmental Engineering. PWS Publishing, 1995.MACINTYRE, A.J. Bombas e Instalações de Bombeamento. LTC, Rio de Janeiro, 1997'
$ws.Range("C22").Value = 'Bibliografia básica:DI BERNARDO, L., Métodos e Técnicas de Tratamento de Água, ABES, Rio de Janeiro, Brasil, 1992.FOUST, A.S., WENZEL, L. A., CLUMP, C.W., MAUS, L., ANDERSEN, L.B. Princípio das operações unitárias. Rio de Janeiro: Editora Guanabara Dois, 1982.GEANKOPLIS, C.J. Procesos de transporte y operaciones unitarias. Compañía Editorial Continental, S.A. de C.V. México, D.F., 1998.PERRY, R.H. and CHILTON, C.H. Manual de Engenharia Química. 5a ed., Guanabara Dois, Rio de Janeiro, 1986.REYNOLDS, T.D.; RICHARDS, P. Unit Operations and Processes in environmental Engineering. PWS Publishing, 1995.MACINTYRE, A.J. Bombas e Instalações de Bombeamento. LTC, Rio de Janeiro, 1997'
